# Changed the dynamic parameters as they come from system properties.
# The APIPATH values that referenced a hard-coded user id / a dependency-test
# generated truid are replaced with a static SYS_USER1 system property, and
# the now-unused helper columns (dependency test id, stored truid, etc.) are
# cleared out. The D2 cell also loses the stray "Hyperlink" look-and-feel
# that was left over from when it used to hold a clickable link.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- APIPATH column (D) -----------------------------------------------
$ws.Range("D2").Value = "/users/user/(SYS_USER1)"
$ws.Range("D3").Value = "/users/user/(SYS_USER1)"
$ws.Range("D4").Value = "/users/user/(SYS_USER1)"
$ws.Range("D5").Value = "/users/user/(SYS_USER1)/image"
$ws.Range("D6").Value = "/users/user/(SYS_USER1)/image"

# D2 used to carry the built-in "Hyperlink" style (blue/underline) since it
# held a clickable-looking GUID path; the new value is a plain string.
$ws.Range("D2").Style = "Normal"

# --- STORE column (K): the test no longer stores a dynamic "truid" ------
$ws.Range("K2").Value = ""

# --- DEPENDENCYTESTS column (I): no longer depends on S1_TC_T1 ----------
$ws.Range("I3").Value = ""
$ws.Range("I4").Value = ""
$ws.Range("I5").Value = ""
$ws.Range("I6").Value = ""

# --- Stray already-blank placeholder cells (no value, no style) get ----
# --- fully dropped from the row once the row is touched again. ---------
$ws.Range("G3").Value = ""
$ws.Range("K3").Value = ""
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = ""
$ws.Range("K4").Value = ""
$ws.Range("G5").Value = ""
$ws.Range("K5").Value = ""
$ws.Range("F6").Value = ""
$ws.Range("G6").Value = ""
$ws.Range("K6").Value = ""

# The now-unused "Hyperlink" built-in cell style is removed entirely.
$wb.Styles.Item("Hyperlink").Delete()

# --- Selection moved from L8 to D2 --------------------------------------
$null = $ws.Range("D2").Select()
